$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "motivation" column (B) used the label "Preventative Health" for three
# rows (Europe, Australia, United States). Per the commit message ("Revised
# Data for consistency"), shorten this label to just "Preventative" so it is
# consistent with the other single-word motivation labels (Wellness, At Risk,
# Sick Role, Self Care).
$oldValue = "Preventative Health"
$newValue = "Preventative"

$used = $ws.UsedRange
foreach ($cell in $used.Cells) {
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
